$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so that numeric-looking
# strings (e.g. "0.638") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.504.78"
$ws.Range("E2").Value = "  +2.04%  "

$ws.Range("D3").Value = "2.366.83"
$ws.Range("E3").Value = "  +6.21%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "317.54"
$ws.Range("E5").Value = "  +7.65%  "

$ws.Range("D6").Value = "107.44"
$ws.Range("E6").Value = "  -4.22%  "

$ws.Range("D7").Value = "0.638"
$ws.Range("E7").Value = "  +2.31%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "0.635"
$ws.Range("E9").Value = "  +4.47%  "

$ws.Range("D10").Value = "42.45"
$ws.Range("E10").Value = "  -4.74%  "

$ws.Range("D11").Value = "0.0934"
$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("D12").Value = "8.68"
$ws.Range("E12").Value = "  -2.45%  "

$ws.Range("E13").Value = "  +2.18%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "16.63"
$ws.Range("E14").Value = "  +9.89%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "0.106"
$ws.Range("E15").Value = "  +2.28%  "

$ws.Range("D16").Value = "2.724.38"
$ws.Range("E16").Value = "  +6.45%  "

$ws.Range("D17").Value = "2.367.19"
$ws.Range("E17").Value = "  +6.36%  "

$ws.Range("D18").Value = "43.507.35"
$ws.Range("E18").Value = "  +2.58%  "

$ws.Range("E19").Value = "  +2.40%  "

$ws.Range("D20").Value = "7.23"
$ws.Range("E20").Value = "  -1.72%  "

$ws.Range("D21").Value = "75.21"
$ws.Range("E21").Value = "  +3.06%  "

$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").Value = "2.52"
$ws.Range("E23").Value = "  +6.34%  "

$ws.Range("D24").Value = "258.64"
$ws.Range("E24").Value = "  +12.59%  "

$ws.Range("D25").Value = "9.29"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").Value = "11.96"
$ws.Range("E26").Value = "  +2.35%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.24"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "38.69"
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").Value = "22.72"
$ws.Range("E30").Value = "  +7.64%  "

$ws.Range("E31").Value = "  -1.37%  "

$ws.Range("D32").Value = "173.51"
$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("D33").Value = "0.0918"
$ws.Range("E33").Value = "  +1.57%  "

$ws.Range("D34").Value = "5.94"
$ws.Range("E34").Value = "  +4.15%  "

$ws.Range("D35").Value = "0.131"
$ws.Range("E35").Value = "  +4.17%  "

$ws.Range("D36").Value = "4.94"
$ws.Range("E36").Value = "  -5.25%  "

$ws.Range("D37").Value = "0.0371"
$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("D38").Value = "4.06"
$ws.Range("E38").Value = "  -6.52%  "

$ws.Range("E39").Value = "  -0.23%  "

$ws.Range("E40").Value = "  +14.32%  "

$ws.Range("D41").Value = "1.50"
$ws.Range("E41").Value = "  +12.91%  "

$ws.Range("D42").Value = "71.47"
$ws.Range("E42").Value = "  -1.41%  "

$ws.Range("E43").Value = "  -1.85%  "

$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").Value = "12.57"
$ws.Range("E45").Value = "  -1.90%  "

$ws.Range("D46").Value = "5.60"
$ws.Range("E46").Value = "  +2.15%  "

$ws.Range("D47").Value = "9.29"
$ws.Range("E47").Value = "  +8.19%  "

$ws.Range("D48").Value = "111.65"
$ws.Range("E48").Value = "  +7.73%  "

$ws.Range("E49").Value = "  -1.43%  "

$ws.Range("E50").Value = "  +2.17%  "

$ws.Range("D51").Value = "0.472"
$ws.Range("E51").Value = "  +6.66%  "
